# Edit slide 4 ("Mechanical Drawing") of the presentation:
#  - "TextBox 13" dimension label changes from 3” to 5”
#  - "TextBox 15" dimension label (2”) together with its two arrow
#    connectors ("Straight Connector 16" and "Straight Connector 17")
#    are removed entirely
#  - "TextBox 18" dimension label changes from 1.5” to 4”

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$s.Shapes.Item("TextBox 13").TextFrame.TextRange.Text = "5”"

$s.Shapes.Item("Straight Connector 16").Delete()
$s.Shapes.Item("Straight Connector 17").Delete()
$s.Shapes.Item("TextBox 15").Delete()

$s.Shapes.Item("TextBox 18").TextFrame.TextRange.Text = "4”"
